$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1297.1111
$ws.Range("I28").Value = 2105.8
$ws.Range("K28").Value = 2105.8
$ws.Range("M28").Value = -1620.8

$ws.Range("H111").Value = 2301.375
$ws.Range("I111").Value = 1076
$ws.Range("J111").Value = 4343.6665
$ws.Range("K111").Value = 3228
$ws.Range("L111").Value = 13030.9995
$ws.Range("M111").Value = -161
$ws.Range("N111").Value = -19164.9995

$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

$ws.Range("H132").Value = 1690.45
$ws.Range("I132").Value = 1253.1052
$ws.Range("K132").Value = 3759.3156
$ws.Range("M132").Value = -1229.3156

$ws.Range("H135").Value = 821.4
$ws.Range("I135").Value = 384.94446
$ws.Range("K135").Value = 3464.50014
$ws.Range("M135").Value = -929.5001400000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10533.958
$ws.Range("I32").Value = 10809.818
$ws.Range("K32").Value = 10809.818
$ws.Range("M32").Value = -10522.818

$ws.Range("H61").Value = 2585.2
$ws.Range("I61").Value = 2585.2
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2585.2
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2373.2
$ws.Range("N61").ClearContents()

$ws.Range("H102").Value = 1079.2
$ws.Range("I102").Value = 1079.2
$ws.Range("K102").Value = 1079.2
$ws.Range("M102").Value = 542.8

$ws.Range("H132").Value = 3748.818
$ws.Range("I132").Value = 2708.1667
$ws.Range("K132").Value = 8124.500100000001
$ws.Range("M132").Value = -5594.500100000001

$ws.Range("H136").Value = 2585.2
$ws.Range("I136").Value = 2585.2
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 7755.599999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5205.599999999999
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 1500
$ws.Range("J22").Value = 1500
$ws.Range("L22").Value = 1500
$ws.Range("N22").Value = -1846

$ws.Range("H86").Value = 2000
$ws.Range("I86").Value = 2000
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 2000
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -877
$ws.Range("N86").ClearContents()

$ws.Range("H89").Value = 2000
$ws.Range("I89").Value = 2000
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 10000
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -4384
$ws.Range("N89").ClearContents()

$ws.Range("H94").Value = 4752.25
$ws.Range("I94").Value = 4836.3335
$ws.Range("J94").Value = 4500
$ws.Range("K94").Value = 4836.3335
$ws.Range("L94").Value = 4500
$ws.Range("M94").Value = -4385.3335
$ws.Range("N94").Value = -5402

$ws.Range("H105").Value = 3950
$ws.Range("I105").Value = 3950
$ws.Range("K105").Value = 3950
$ws.Range("M105").Value = -2203

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4000
$ws.Range("I16").Value = 4000
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 4000
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -3713
$ws.Range("N16").ClearContents()

$ws.Range("H31").Value = 1851.2858
$ws.Range("J31").Value = 2124.5
$ws.Range("L31").Value = 2124.5
$ws.Range("N31").Value = -2714.5

$ws.Range("H34").Value = 1851.2858
$ws.Range("J34").Value = 2124.5
$ws.Range("L34").Value = 2124.5
$ws.Range("N34").Value = -2528.5

$ws.Range("H113").Value = 4000
$ws.Range("I113").Value = 4000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -1830
$ws.Range("N113").ClearContents()

$ws.Range("H132").Value = 2202.6
$ws.Range("J132").Value = 2862.8572
$ws.Range("L132").Value = 8588.5716
$ws.Range("N132").Value = -13648.5716

$ws.Range("H134").Value = 3917.5
$ws.Range("I134").Value = 3917.5
$ws.Range("K134").Value = 11752.5
$ws.Range("M134").Value = -9217.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 1026.909
$ws.Range("I86").Value = 1327.25
$ws.Range("J86").Value = 855.2857
$ws.Range("K86").Value = 3981.75
$ws.Range("L86").Value = 2565.8571
$ws.Range("M86").Value = -2795.75
$ws.Range("N86").Value = -4937.8571

$ws.Range("H87").Value = 5002.5557
$ws.Range("I87").Value = 3858.8572
$ws.Range("K87").Value = 11576.5716
$ws.Range("M87").Value = -10328.5716

$ws.Range("H89").Value = 1026.909
$ws.Range("I89").Value = 1327.25
$ws.Range("J89").Value = 855.2857
$ws.Range("K89").Value = 11945.25
$ws.Range("L89").Value = 7697.571300000001
$ws.Range("M89").Value = -6017.25
$ws.Range("N89").Value = -19553.5713

$ws.Range("H90").Value = 5002.5557
$ws.Range("I90").Value = 3858.8572
$ws.Range("K90").Value = 34729.7148
$ws.Range("M90").Value = -28489.7148

$ws.Range("H121").Value = 196.125
$ws.Range("I121").Value = 231
$ws.Range("K121").Value = 693
$ws.Range("M121").Value = 617

$ws.Range("H129").Value = 1933.3
$ws.Range("J129").Value = 2041.625
$ws.Range("L129").Value = 6124.875
$ws.Range("N129").Value = -16124.875

$ws.Range("H131").Value = 2664.5386
$ws.Range("J131").Value = 2629.9
$ws.Range("L131").Value = 7889.700000000001
$ws.Range("N131").Value = -17969.7

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 399.875
$ws.Range("I2").Value = 500
$ws.Range("J2").Value = 339.8
$ws.Range("K2").Value = 500
$ws.Range("L2").Value = 339.8
$ws.Range("M2").Value = -387
$ws.Range("N2").Value = -565.8

$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws.Range("H97").Value = 1264
$ws.Range("I97").Value = 1341.1666
$ws.Range("K97").Value = 1341.1666
$ws.Range("M97").Value = -845.1666

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1834.909
$ws.Range("I93").Value = 1107.7142
$ws.Range("K93").Value = 1107.7142
$ws.Range("M93").Value = 140.2858000000001

$ws.Range("H95").Value = 39990
$ws.Range("J95").Value = 39990
$ws.Range("L95").Value = 39990
$ws.Range("N95").Value = -45482

$ws.Range("H100").Value = 5183.3335
$ws.Range("I100").Value = 5183.3335
$ws.Range("K100").Value = 5183.3335
$ws.Range("M100").Value = -4642.3335

$ws.Range("H132").Value = 5391.231
$ws.Range("I132").Value = 3261.25
$ws.Range("J132").Value = 8799.2
$ws.Range("K132").Value = 9783.75
$ws.Range("L132").Value = 26397.6
$ws.Range("M132").Value = -7253.75
$ws.Range("N132").Value = -31457.6

$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 5562.3335
$ws.Range("I136").Value = 5562.3335
$ws.Range("K136").Value = 16687.0005
$ws.Range("M136").Value = -14137.0005
